# Burndown chart / review report update
# - Corrects the actual ("effettivo") story points recorded for Sprint 1-3
#   on the "Burndown Chart" sheet (planned vs effective burndown columns).
#   The "Burndown ideale"/"Burndown effettivo" (D/E) columns are formulas
#   that recompute automatically once the inputs below are written.
# - Leaves the "Burndown Chart" tab selected/active (instead of
#   "Task Sprint 1"), with an updated zoom level and cell selection,
#   mirroring where the author was last looking in the workbook.

$wb = $excel.ActiveWorkbook

$wsBurndown = $wb.Worksheets.Item("Burndown Chart")

# --- Update the effective story-points data (B/C columns) -----------------
# B5 (Sprint 3 planned):            28  -> 0
$wsBurndown.Range("B5").Value = 0
# C3 (Sprint 1 effective):          20  -> 22
$wsBurndown.Range("C3").Value = 22
# C4 (Sprint 2 effective):          25  -> 0
$wsBurndown.Range("C4").Value = 0
# C5 (Sprint 3 effective):          30  -> 0
$wsBurndown.Range("C5").Value = 0

# D2:D5 / E2:E5 are formulas ( J2 / prior - current ) and recalculate on
# their own; no need (and no way) to poke their cached values directly.

# --- Switch the active sheet / selection back to "Burndown Chart" ---------
[void]$wsBurndown.Activate()
$excel.ActiveWindow.Zoom = 98
[void]$wsBurndown.Range("C5").Select()
